# Update common mock files for UT IAC JR regional lists
#
# Rebuilds Sheet1 with the new "Venue" based hearing-list layout:
# Venue | Judge(s) | Hearing time | Case reference number | Case title | Hearing type | Additional information

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe all existing content + column formatting (A:H) so we start from a
# clean sheet, matching the completely reshaped table in the target file.
$ws.Columns("A:H").Delete()

# Header row
$ws.Range("A1").Value = "Venue"
$ws.Range("B1").Value = "Judge(s)"
$ws.Range("C1").Value = "Hearing time"
$ws.Range("D1").Value = "Case reference number"
$ws.Range("E1").Value = "Case title"
$ws.Range("F1").Value = "Hearing type"
$ws.Range("G1").Value = "Additional information"

# Row 2
$ws.Range("A2").Value = "Venue A"
$ws.Range("B2").Value = "Judge A"
$ws.Range("C2").Value = "10:30am"
$ws.Range("D2").Value = 1234
$ws.Range("E2").Value = "Case title A"
$ws.Range("F2").Value = "Hearing type A"
$ws.Range("G2").Value = "This is additional information"

# Row 3
$ws.Range("A3").Value = "Venue B"
$ws.Range("B3").Value = "Judge B"
$ws.Range("C3").Value = "11am"
$ws.Range("D3").Value = 4567
$ws.Range("E3").Value = "Case title B"
$ws.Range("F3").Value = "Hearing type B"
$ws.Range("G3").Value = "This is additional information"

# Row 4
$ws.Range("A4").Value = "Venue C"
$ws.Range("B4").Value = "Judge C"
$ws.Range("C4").Value = "11:30am"
$ws.Range("D4").Value = 5678
$ws.Range("E4").Value = "Case title C"
$ws.Range("F4").Value = "Hearing type C"
$ws.Range("G4").Value = "This is additional information"

# The hearing-time column is formatted (text stored, but tagged with a
# time number format), matching numFmtId 20 ("h:mm") in the target styles.
$ws.Range("C2:C4").NumberFormat = "h:mm"

# Column widths for the new layout (best-fit-like widths from the source file)
$ws.Columns("C").ColumnWidth = 10.58333333333334
$ws.Columns("D").ColumnWidth = 19.250000000000007
$ws.Columns("E").ColumnWidth = 12.583333333333343
$ws.Columns("F").ColumnWidth = 16.750000000000014
$ws.Columns("G").ColumnWidth = 18.250000000000007

# Match the saved selection in the target file
[void]$ws.Range("D8").Select()
